# Personal History Statement - edits per commit "update. d-1 day. santa cruz"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Font rename: "times" -> "Times" for every run (ascii + hAnsi), leaving
#    eastAsia ("Arial Unicode MS") and cs ("Times New Roman") untouched.
#    Setting Font.Name on each paragraph's range rewrites both ascii/hAnsi
#    consistently (including the paragraph-mark rPr in pPr).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs.Item($i).Range.Font.Name = "Times"
}

# ---------------------------------------------------------------------------
# 2) Paragraph 1: the trailing "_GoBack" bookmark around the lone space is
#    gone and the three runs it split become one. Re-saving the identical
#    text over the whole span collapses the runs (and drops the bookmark
#    that is inside the replaced range).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "My interests in computer arouse when I was a little boy no more than 10 years old. At that time a computer was merely an entertainment for me and video games are all the things it can offer to me. Games, in my very personal opinion at that time, are entertaining but not that useful and hence computer did not have a greater meaning. Until years later, one day I accidently opened a map editor of a game (Star Craft) and saw a complex tool with multiple leveled menu and hundreds of buttons that can literally control every piece on the screen. Through that map editor I imagined a whole different usage of a computer, that is, it can run powerful software to perform and organize complicated tasks far beyond games. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "My interests in computer arouse when I was a little boy no more than 10 years old. At that time a computer was merely an entertainment for me and video games are all the things it can offer to me. Games, in my very personal opinion at that time, are entertaining but not that useful and hence computer did not have a greater meaning. Until years later, one day I accidently opened a map editor of a game (Star Craft) and saw a complex tool with multiple leveled menu and hundreds of buttons that can literally control every piece on the screen. Through that map editor I imagined a whole different usage of a computer, that is, it can run powerful software to perform and organize complicated tasks far beyond games. ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Paragraph 2: the "_GoBack" bookmark re-appears around "high school
#    graduation" (this is where the author's cursor ended up). Bookmarks.Add
#    with the reserved name "_GoBack" replaces the old one wherever it was.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("high school graduation", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# ---------------------------------------------------------------------------
# 4) Paragraph 3 (ethnic identity): insert "the beginning year of " before
#    "my university life" and pluralize "issue" -> "issues".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("my university life also presented", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($rng.Start, $rng.Start)
$insPoint.InsertBefore("the beginning year of ")

$rng = $d.Content
$rng.Find.Execute("ethnic identity issue", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter("s")

# ---------------------------------------------------------------------------
# 5) Paragraph 5 (parents' hardships): fix the typo "go though" -> "go
#    through", then add a comma after "...ceived college education".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("had to go though.", $false, $false, $false, $false, $false, $true, 1, $false, "had to go through.", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("ceived college education", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter(",")

# ---------------------------------------------------------------------------
# 6) Final paragraph: "I would be proud since it will be a right one." ->
#    "I would be proud of myself since it would prove to be a right one."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("would be proud since it will be a right one.", $false, $false, $false, $false, $false, $true, 1, $false, "would be proud of myself since it would prove to be a right one.", 2) | Out-Null
